$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for new rows from existing same-style rows ---
$xlPasteFormats = -4122

# Style "2" rows (84-88) <- copy format from row 77 (A:D)
$ws.Range("A77:D77").Copy()
$ws.Range("A84:D88").PasteSpecial($xlPasteFormats)

# Style "1" rows (89-96) <- copy format from row 79 (A:D)
$ws.Range("A79:D79").Copy()
$ws.Range("A89:D96").PasteSpecial($xlPasteFormats)

# Style "4" (note cells in column E) <- copy format from E61
$ws.Range("E61").Copy()
$ws.Range("E85").PasteSpecial($xlPasteFormats)
$ws.Range("E90").PasteSpecial($xlPasteFormats)
$ws.Range("E96").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Set cell values in the precise order that reproduces the original
#     shared-string insertion sequence ---
$ws.Range("B84").Value = 'dashboard/ver-publicacion-ofrecida/idPublicacion'
$ws.Range("C84").Value = 'Muestra todos los datos de la publicación.'
$ws.Range("C85").Value = 'Muestra todos los datos del cliente dueño de la publicación. Muetra puntaje de la publicación y del servicio en general'
$ws.Range("A84").Value = 'Listado datos de la publicación'
$ws.Range("A85").Value = 'Listado datos del cliente dueño de la publicación'
$ws.Range("A86").Value = 'Listado de los comentarios y puntuaciones de la publicación'
$ws.Range("C86").Value = 'Muestra todos los comentarios y puntajes de la publicación'
$ws.Range("A87").Value = 'Responder comentario realizado a la publicación'
$ws.Range("C87").Value = 'Click en link responder. No se ingresa comentario. Solo se responde si es el dueño de la publicación.'
$ws.Range("C88").Value = 'Click en link responder. Comentario correcto. Solo se responde si es el dueño de la publicación.'
$ws.Range("B89").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("A89").Value = 'Listado datos de la solicitud'
$ws.Range("C89").Value = 'Muestra todos los datos de la solicitud.'
$ws.Range("C90").Value = 'Muestra imagen y nombre de usuario del cliente dueño de la publicación.'
$ws.Range("E90").Value = 'FALTA MOSTRAR UBICACIÓN'
$ws.Range("E85").Value = 'VER LOS PUNTAJES, EL DE LOS SERVICIOS DEBE SER DE SOLICITUDES Y OFERTAS. FALTA MOSTRAR UBICACIÓN'
$ws.Range("A91").Value = 'Listado de todas las propuestas'
$ws.Range("C91").Value = 'Muestra todos las propuestas realizadas hasta el momento.'
$ws.Range("A92").Value = 'Realizar una postulación'
$ws.Range("C92").Value = 'Click en link postularme. No se ingresa datos.'
$ws.Range("D92").Value = 'Alert correspondiente indicando que se debe ingresar un texto.'
$ws.Range("C93").Value = 'Click en link postularme. Datos correctos.'
$ws.Range("A95").Value = 'Mostrar/Ocultar Propuestas'
$ws.Range("C95").Value = 'Muestra/Oculta las propuestas realizadas.'
$ws.Range("A96").Value = 'Aceptar propuesta'
$ws.Range("C96").Value = 'Click en aceptar propuesta. Unicamente el due;o de la solicitud tiene el link.'
$ws.Range("D96").Value = 'OK. Se acepta la propuesta, se finaliza la solicitud, se habilita la calificaci''on del usuario contratado.'
$ws.Range("E96").Value = 'FALTA NOTIFICAR AL CONTRATADO, MOSTRAR DATOS DEL DUE;O DE LA PUBLICACI''ON. MOSTRAR DATOS DEL TRABAJADOR.'
$ws.Range("D84").Value = 'OK.'
$ws.Range("B85").Value = 'dashboard/ver-publicacion-ofrecida/idPublicacion'
$ws.Range("D85").Value = 'OK.'
$ws.Range("B86").Value = 'dashboard/ver-publicacion-ofrecida/idPublicacion'
$ws.Range("D86").Value = 'OK.'
$ws.Range("B87").Value = 'dashboard/ver-publicacion-ofrecida/idPublicacion'
$ws.Range("D87").Value = 'Alert correspondiente indicando que se debe ingresar un comentario.'
$ws.Range("A88").Value = 'Responder comentario realizado a la publicación'
$ws.Range("B88").Value = 'dashboard/ver-publicacion-ofrecida/idPublicacion'
$ws.Range("D88").Value = 'OK.'
$ws.Range("D89").Value = 'OK.'
$ws.Range("A90").Value = 'Listado datos del cliente dueño de la publicación'
$ws.Range("B90").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("D90").Value = 'OK.'
$ws.Range("B91").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("D91").Value = 'OK.'
$ws.Range("B92").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("A93").Value = 'Realizar una postulación'
$ws.Range("B93").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("D93").Value = 'OK.'
$ws.Range("A94").Value = 'Responder comentario realizado a la publicación'
$ws.Range("B94").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("C94").Value = 'Click en link responder. Comentario correcto. Solo se responde si es el dueño de la publicación.'
$ws.Range("D94").Value = 'OK.'
$ws.Range("B95").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'
$ws.Range("D95").Value = 'OK.'
$ws.Range("B96").Value = 'dashboard/ver-publicacion-solicitada/idPublicacion'

# --- Column widths (best-effort; runtime rounds to pixel-based increments) ---
$ws.Columns.Item(3).ColumnWidth = 101.92
$ws.Columns.Item(5).ColumnWidth = 107.59

# --- View / selection update ---
$win = $excel.ActiveWindow
$win.ScrollRow = 64
$win.ScrollColumn = 1
$ws.Range("E96").Select()
